$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data: Numero, Nombre, Horas (formatted 0.00), Horas extra
$data = @(
    @(12, "Pedro Perez",       80.5, 4),
    @(32, "Juan Gopnzalez",    80.5, 2.3),
    @(54, "Pablo Picapiedras", 80.5, 5.1),
    @(85, "Marcelo Gomez",     80.5, 3.4),
    @(64, "Alberto Paredes",   80.5, 2.5),
    @(87, "Sebastian Romani",  80.5, 5),
    @(15, "Mauricio Morales",  80.5, 81),
    @(28, "Facundo Fredes",    80.5, 3.5)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $rowVals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowVals[0]
    $ws.Cells.Item($row, 2).Value = $rowVals[1]
    $ws.Cells.Item($row, 3).Value = $rowVals[2]
    $ws.Cells.Item($row, 4).Value = $rowVals[3]
}

# Column C (Horas) carries the "0.00" number format
$ws.Range("C1:C8").NumberFormat = "0.00"

# Column widths matching the authored layout
$ws.Columns.Item(1).ColumnWidth = 6
$ws.Columns.Item(2).ColumnWidth = 16

# Restore the selection left by the author (cell below the data)
$ws.Range("D9").Select() | Out-Null

Write-Output "Done"
